# Applies the 2002_bacias.xlsx changes:
#  - header row retitled (LEGENDA -> legenda, area_km2 -> area, area_km2_1 -> area_km2)
#  - column B (area_km2) data values updated for rows 2-26

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates ---
$ws.Cells.Item(1, 1).Value = "legenda"
$ws.Cells.Item(1, 2).Value = "area"
$ws.Cells.Item(1, 3).Value = "nome"
$ws.Cells.Item(1, 4).Value = "area_km2"

# --- Updated area_km2 values for rows 2-26 ---
$newValues = @(
    13141.9972595,
    1638.95622443,
    9254.16171808,
    1567.77845243,
    2116.88353751,
    3942.25013061,
    7150.98928162,
    6328.53226945,
    1849.3153882,
    884.441544154,
    10245.6747183,
    2321.51552537,
    328.721732755,
    1997.06759274,
    17473.66252,
    4106.52340143,
    10718.4168586,
    5301.42115388,
    9874.269802430001,
    9820.44094998,
    4166.36667473,
    9852.6254023,
    8445.983877320001,
    1337.14824275,
    14911.9201979
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}
